# correção nos dados e inicio da analise PNAD 2009
#
# The sheet originally had two "header-only" rows (row 5 = "situação do
# domicílio", row 8 = "grandes regiões e unidades da federação") that were
# category sub-headers with no data of their own, followed immediately by
# the real data rows below them. The fix removes those two label-only rows
# (shifting every row below them up), which both corrects the data
# alignment (each region label now lines up with its own row of numbers)
# and drops the two now-unused caption strings from the shared string
# table. It also relabels the "total" column header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("situação do domicílio") was a label-only row directly above the
# "urbana"/"rural" data rows below it - remove it so "urbana" takes row 5.
$ws.Rows("5").Delete()

# After the row-5 deletion, row 8 is no longer "rondônia" - the old row 8
# ("grandes regiões e unidades da federação", also label-only) is now row 7.
# Remove it too so the region data shifts up into place.
$ws.Rows("7").Delete()

# The "unnamed: 1_level_1" placeholder header becomes "total".
$ws.Range("B2").Value2 = "total"
